$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per diff (column D values are numeric-looking and must be
# written as exact text so Excel does not normalize the literal representation,
# e.g. "3.700" -> 3.7 or "0.0001500" -> 0.00015)
$textUpdates = [ordered]@{
    'D2' = '250.34'
    'D3' = '22.84'
    'D4' = '5.483'
    'D5' = '0.05658'
    'D6' = '3.423'
    'D7' = '6.369'
    'D8' = '0.8159'
    'D9' = '0.9420'
    'D10' = '0.1442'
    'D11' = '0.07539'
    'D12' = '0.03117'
    'D13' = '0.03099'
    'D14' = '0.09366'
    'D15' = '3.560'
    'D16' = '0.001582'
    'D17' = '0.04767'
    'D18' = '0.006401'
    'D19' = '0.004997'
    'D20' = '0.001031'
    'D21' = '0.0001500'
    'D22' = '3.700'
    'D23' = '2.192'
    'D24' = '0.01160'
    'D26' = '0.1295'
    'D28' = '0.0003031'
    'D40' = '0.04041'
    'D41' = '0.006758'
    'D42' = '0.1070'
    'D43' = '0.002719'
    'D44' = '0.007558'
    'D45' = '0.00005801'
    'D47' = '0.4999'
    'D49' = '0.00002100'
}

$plainUpdates = [ordered]@{
    'B18' = 'TigerCash'
    'C18' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'E18' = '17TigerCashTCH'
    'B19' = 'HotbitToken'
    'C19' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'E19' = '18HotbitTokenHTB'
    'B20' = 'BitKan'
    'C20' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'E20' = '19BitKanKAN'
    'B21' = 'NitroEx'
    'C21' = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
    'E21' = '20NitroExNTX'
    'B22' = 'LEO'
    'C22' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'E22' = '21LEOLEO'
    'B23' = 'BTSEToken'
    'C23' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'E23' = '22BTSETokenBTSE'
    'B24' = 'One'
    'C24' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'E24' = '23OneONEBestin24h'
    'E48' = '47BOLOBOLO'
}

foreach ($ref in $textUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$ref]
    $cell.Style = "Normal"
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}
